$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C ("Description") to hold a new
# "Term Type" field. This shifts the previous C:Q columns to D:R.
$ws.Columns("C").Insert()

# Header for the newly inserted column.
$ws.Range("C1").Value = "Term Type"

# Populate the new column's data rows with the term type values.
$ws.Range("C2").Value = "phenotype"
$ws.Range("C3").Value = "germplasm attribute"
$ws.Range("C6").Value = "Germplasm Passport"

# Match the bordered look of the rest of the data table for the new
# column (rows 4 and 5 have no Term Type value, but still get the
# border formatting applied down the column) by copying the format
# from the neighboring Description column, then removing its fill so
# the new column stays unshaded.
$ws.Range("D2:D6").Copy()
$ws.Range("C2:C6").PasteSpecial(-4122)
$ws.Range("C2:C6").Interior.Pattern = -4142
